# Scheduled runner: refresh Universalis market-price snapshots
# for the per-class Leve profit tables (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) in place.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1715
$ws.Range("I40").Value = 1595.2727
$ws.Range("J40").Value = 1834.7273
$ws.Range("K40").Value = 1595.2727
$ws.Range("L40").Value = 1834.7273
$ws.Range("M40").Value = -1420.2727
$ws.Range("N40").Value = -2184.7273
$ws.Range("H74").Value = 2652.3235
$ws.Range("I74").Value = 2192.7334
$ws.Range("K74").Value = 2192.7334
$ws.Range("M74").Value = -1256.7334
$ws.Range("H77").Value = 2652.3235
$ws.Range("I77").Value = 2192.7334
$ws.Range("K77").Value = 10963.667
$ws.Range("M77").Value = -6283.667000000001
$ws.Range("H133").Value = 51528.57
$ws.Range("J133").Value = 51528.57
$ws.Range("L133").Value = 51528.57
$ws.Range("N133").Value = -61648.57
$ws.Range("H137").Value = 3429.24
$ws.Range("I137").Value = 984.65
$ws.Range("J137").Value = 5058.967
$ws.Range("K137").Value = 2953.95
$ws.Range("L137").Value = 15176.901
$ws.Range("M137").Value = -403.9499999999998
$ws.Range("N137").Value = -20276.901

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1881.2059
$ws.Range("I45").Value = 1509.3529
$ws.Range("K45").Value = 1509.3529
$ws.Range("M45").Value = -1132.3529
$ws.Range("H104").Value = 41056
$ws.Range("J104").Value = 41056
$ws.Range("L104").Value = 41056
$ws.Range("N104").Value = -48044
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("H130").Value = 33186.555
$ws.Range("J130").Value = 33186.555
$ws.Range("L130").Value = 33186.555
$ws.Range("N130").Value = -43226.555
$ws.Range("N121").ClearContents()

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 200
$ws.Range("I22").Value = 200
$ws.Range("K22").Value = 200
$ws.Range("M22").Value = -27
$ws.Range("H96").Value = 20747.75
$ws.Range("I96").Value = 7916.6665
$ws.Range("J96").Value = 28446.4
$ws.Range("K96").Value = 7916.6665
$ws.Range("L96").Value = 28446.4
$ws.Range("M96").Value = -5170.6665
$ws.Range("N96").Value = -33938.4

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 41759.668
$ws.Range("J20").Value = 41759.668
$ws.Range("L20").Value = 41759.668
$ws.Range("N20").Value = -42231.668
$ws.Range("H30").Value = 41759.668
$ws.Range("J30").Value = 41759.668
$ws.Range("L30").Value = 41759.668
$ws.Range("N30").Value = -41941.668
$ws.Range("H31").Value = 3882.0354
$ws.Range("I31").Value = 1504.65
$ws.Range("J31").Value = 4613.5386
$ws.Range("K31").Value = 1504.65
$ws.Range("L31").Value = 4613.5386
$ws.Range("M31").Value = -1209.65
$ws.Range("N31").Value = -5203.5386
$ws.Range("H34").Value = 3882.0354
$ws.Range("I34").Value = 1504.65
$ws.Range("J34").Value = 4613.5386
$ws.Range("K34").Value = 1504.65
$ws.Range("L34").Value = 4613.5386
$ws.Range("M34").Value = -1302.65
$ws.Range("N34").Value = -5017.5386
$ws.Range("H58").Value = 1183.8572
$ws.Range("I58").Value = 463.63635
$ws.Range("J58").Value = 1649.8823
$ws.Range("K58").Value = 463.63635
$ws.Range("L58").Value = 1649.8823
$ws.Range("M58").Value = -260.63635
$ws.Range("N58").Value = -2055.8823
$ws.Range("H62").Value = 4229.8125
$ws.Range("I62").Value = 4694.3335
$ws.Range("J62").Value = 2836.25
$ws.Range("K62").Value = 4694.3335
$ws.Range("L62").Value = 2836.25
$ws.Range("M62").Value = -4070.3335
$ws.Range("N62").Value = -4084.25
$ws.Range("H65").Value = 4229.8125
$ws.Range("I65").Value = 4694.3335
$ws.Range("J65").Value = 2836.25
$ws.Range("K65").Value = 23471.6675
$ws.Range("L65").Value = 14181.25
$ws.Range("M65").Value = -20351.6675
$ws.Range("N65").Value = -20421.25
$ws.Range("H128").Value = 41759.668
$ws.Range("J128").Value = 41759.668
$ws.Range("L128").Value = 41759.668
$ws.Range("N128").Value = -51719.668
$ws.Range("H132").Value = 40006390
$ws.Range("I132").Value = 52639230
$ws.Range("J132").Value = 2392.6667
$ws.Range("K132").Value = 157917690
$ws.Range("L132").Value = 7178.000100000001
$ws.Range("M132").Value = -157915160
$ws.Range("N132").Value = -12238.0001
$ws.Range("H134").Value = 2318
$ws.Range("I134").Value = 2899.8333
$ws.Range("J134").Value = 1819.2858
$ws.Range("K134").Value = 8699.499899999999
$ws.Range("L134").Value = 5457.857400000001
$ws.Range("M134").Value = -6164.499899999999
$ws.Range("N134").Value = -10527.8574
$ws.Range("H136").Value = 1183.8572
$ws.Range("I136").Value = 463.63635
$ws.Range("J136").Value = 1649.8823
$ws.Range("K136").Value = 1390.90905
$ws.Range("L136").Value = 4949.6469
$ws.Range("M136").Value = 1159.09095
$ws.Range("N136").Value = -10049.6469

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 2540
$ws.Range("J75").Value = 2540
$ws.Range("L75").Value = 7620
$ws.Range("N75").Value = -9616
$ws.Range("H78").Value = 2540
$ws.Range("J78").Value = 2540
$ws.Range("L78").Value = 22860
$ws.Range("N78").Value = -32844
$ws.Range("H131").Value = 688.8511
$ws.Range("J131").Value = 913.2
$ws.Range("L131").Value = 2739.6
$ws.Range("N131").Value = -12819.6

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 26666.666
$ws.Range("I64").Value = 10000
$ws.Range("K64").Value = 10000
$ws.Range("M64").Value = -9752
$ws.Range("H67").Value = 26666.666
$ws.Range("I67").Value = 10000
$ws.Range("K67").Value = 10000
$ws.Range("M67").Value = -9142
$ws.Range("H107").Value = 350
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 350
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 350
$ws.Range("N107").Value = -4190
$ws.Range("H126").Value = 1886.8334
$ws.Range("I126").Value = 1350.5
$ws.Range("J126").Value = 2155
$ws.Range("K126").Value = 4051.5
$ws.Range("L126").Value = 6465
$ws.Range("M126").Value = -1581.5
$ws.Range("N126").Value = -11405
$ws.Range("M107").ClearContents()

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1783.3077
$ws.Range("I7").Value = 1632.8
$ws.Range("J7").Value = 2285
$ws.Range("K7").Value = 1632.8
$ws.Range("L7").Value = 2285
$ws.Range("M7").Value = -1520.8
$ws.Range("N7").Value = -2509
$ws.Range("H40").Value = 3373.9375
$ws.Range("I40").Value = 3360
$ws.Range("J40").Value = 3380.2727
$ws.Range("K40").Value = 3360
$ws.Range("L40").Value = 3380.2727
$ws.Range("M40").Value = -3224
$ws.Range("N40").Value = -3652.2727
$ws.Range("H46").Value = 2655.6365
$ws.Range("J46").Value = 2828.2
$ws.Range("L46").Value = 2828.2
$ws.Range("N46").Value = -3204.2
$ws.Range("H126").Value = 1783.3077
$ws.Range("I126").Value = 1632.8
$ws.Range("J126").Value = 2285
$ws.Range("K126").Value = 4898.4
$ws.Range("L126").Value = 6855
$ws.Range("M126").Value = -2428.4
$ws.Range("N126").Value = -11795
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 40001080
$ws.Range("H132").Value = 5148.3228
$ws.Range("I132").Value = 5372.1724
$ws.Range("J132").Value = 1902.5
$ws.Range("K132").Value = 16116.5172
$ws.Range("L132").Value = 5707.5
$ws.Range("M132").Value = -13586.5172
$ws.Range("N132").Value = -10767.5
$ws.Range("H136").Value = 2630.5796
$ws.Range("I136").Value = 2578.3455
$ws.Range("J136").Value = 2835.7856
$ws.Range("K136").Value = 7735.0365
$ws.Range("L136").Value = 8507.356800000001
$ws.Range("M136").Value = -5185.0365
$ws.Range("N136").Value = -13607.3568

